# Applies weekly fruit/vegetable (Pomelo) price update:
# the data rows (2-9) get their date/volume/price values permuted to
# reflect the refreshed weekly consolidated figures, while the
# "Unidad de comercializacion" (Q) text is updated independently per row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns whose values move together as one "record" for each source row,
# keeping A,B,C,E,F,G,H,I,J,K,L,R,T fixed per row. Q is handled separately
# below since it does not travel with the rest of the record.
$cols = @("D", "M", "N", "O", "P", "S")

# Capture the current (pre-edit) values for each affected row/column.
$data = @{}
foreach ($row in 2..9) {
    $rowData = @{}
    foreach ($col in $cols) {
        $rowData[$col] = $ws.Range("$col$row").Value2
    }
    $data[$row] = $rowData
}

# Mapping of source row -> destination row (permutation of the 8 records).
$mapping = @{
    2 = 3
    3 = 2
    4 = 8
    5 = 9
    6 = 5
    7 = 7
    8 = 6
    9 = 4
}

foreach ($srcRow in $mapping.Keys) {
    $dstRow = $mapping[$srcRow]
    $rowData = $data[$srcRow]
    foreach ($col in $cols) {
        $ws.Range("$col$dstRow").Value2 = $rowData[$col]
    }
}

# Final "Unidad de comercializacion" text per row.
$qValues = @{
    2 = "`$/caja 14 kilos empedrada"
    3 = "`$/caja 14 kilos empedrada"
    4 = "`$/caja 14 kilos"
    5 = "`$/caja 14 kilos empedrada"
    6 = "`$/caja 14 kilos"
    7 = "`$/caja 14 kilos empedrada"
    8 = "`$/caja 14 kilos empedrada"
    9 = "`$/caja 14 kilos empedrada"
}

foreach ($row in $qValues.Keys) {
    $ws.Range("Q$row").Value2 = $qValues[$row]
}
